# Update "Forecast Comparison" sheet:
#  - insert a new "Week_Start_Date" column after "Week" (shifts ASIN.. right by one)
#  - change Week labels from zero-padded ("W01") to unpadded ("W1")
#  - fill the new Week_Start_Date column with the week's start date (as text)
#  - change is_holiday_week values from numeric 0 to boolean FALSE

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column at B; existing B (ASIN) through I (is_holiday_week) shift to C..J
$ws.Columns.Item(2).Insert()

# New header
$ws.Range("B1").Value = "Week_Start_Date"

# Week number (column A) and week-start-date (column B) values for rows 2-17
$weeks = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16)
$dates = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)

# Make sure column B is treated as plain text so the dates are stored as strings
$colB = $ws.Range("B2:B17")
$colB.NumberFormat = "@"

for ($i = 0; $i -lt $weeks.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = "W" + $weeks[$i]
    $ws.Cells.Item($row, 2).Value = $dates[$i]
}

# Drop the text-format marker so the cells look like ordinary (unstyled) text cells
$colB.ClearFormats()

# is_holiday_week is now column J (after the insert) - convert 0 -> FALSE (boolean)
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 10).Value = $false
}
